$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 20.41270024389399
$ws.Range("C2").Value = 17.26868271090615
$ws.Range("D2").Value = 6.953944436805044
$ws.Range("E2").Value = 12.71111770342536
$ws.Range("F2").Value = 43.92638751237448
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 27.5058778154737
$ws.Range("J2").Value = 10.28437714029083
$ws.Range("N2").Value = 18.53984072357984

# Row 3
$ws.Range("B3").Value = 19.90059342218844
$ws.Range("C3").Value = 16.76006302165466
$ws.Range("D3").Value = 6.94784728612303
$ws.Range("E3").Value = 12.68434590529637
$ws.Range("F3").Value = 43.72777172053112
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 27.52309948137522
$ws.Range("J3").Value = 10.29380399114021
$ws.Range("N3").Value = 18.61783209767465

# Row 4
$ws.Range("B4").Value = 19.5854085791532
$ws.Range("C4").Value = 16.44545893365876
$ws.Range("D4").Value = 6.944803949421152
$ws.Range("E4").Value = 12.6705795576984
$ws.Range("F4").Value = 43.61938210404897
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 27.54139620639288
$ws.Range("J4").Value = 10.30152207838112
$ws.Range("N4").Value = 18.66772759342791

# Row 5
$ws.Range("B5").Value = 19.45698718252576
$ws.Range("C5").Value = 16.31689344150379
$ws.Range("D5").Value = 6.943741062017095
$ws.Range("E5").Value = 12.66564429437032
$ws.Range("F5").Value = 43.57864446926092
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 27.55078441393816
$ws.Range("J5").Value = 10.30515175345569
$ws.Range("N5").Value = 18.68856731357358

# Row 6
$ws.Range("B6").Value = 19.43567025350202
$ws.Range("C6").Value = 16.29553001385224
$ws.Range("D6").Value = 6.943575312388017
$ws.Range("E6").Value = 12.66486561524079
$ws.Range("F6").Value = 43.57208783411529
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 27.55245970965654
$ws.Range("J6").Value = 10.30578369402346
$ws.Range("N6").Value = 18.69205840309167

# Row 7
$ws.Range("B7").Value = 19.58367628724852
$ws.Range("C7").Value = 16.4437262213625
$ws.Range("D7").Value = 6.944788895539627
$ws.Range("E7").Value = 12.67051026412818
$ws.Range("F7").Value = 43.61881878141001
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 27.54151500977559
$ws.Range("J7").Value = 10.30156906894216
$ws.Range("N7").Value = 18.66800659039925

# Row 8
$ws.Range("B8").Value = 20.23640336403734
$ws.Range("C8").Value = 17.09391766427245
$ws.Range("D8").Value = 6.951697513356546
$ws.Range("E8").Value = 12.70133415287311
$ws.Range("F8").Value = 43.85510493491164
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 27.51020747147289
$ws.Range("J8").Value = 10.28722650177254
$ws.Range("N8").Value = 18.56631641237104

# Row 9
$ws.Range("B9").Value = 21.50176199423407
$ws.Range("C9").Value = 18.34163096372701
$ws.Range("D9").Value = 6.970748427323466
$ws.Range("E9").Value = 12.78282992754745
$ws.Range("F9").Value = 44.42484994014125
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 27.51047921114985
$ws.Range("J9").Value = 10.27444637018347
$ws.Range("N9").Value = 18.38275217405637

# Row 10
$ws.Range("B10").Value = 22.41155779875934
$ws.Range("C10").Value = 19.23062506491059
$ws.Range("D10").Value = 6.988025867661277
$ws.Range("E10").Value = 12.85530515753779
$ws.Range("F10").Value = 44.90632413087205
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 27.54873848264491
$ws.Range("J10").Value = 10.27444983223694
$ws.Range("N10").Value = 18.25742626349261

# Row 11
$ws.Range("B11").Value = 22.81904576916434
$ws.Range("C11").Value = 19.62699384937129
$ws.Range("D11").Value = 6.996581307704076
$ws.Range("E11").Value = 12.89095088549031
$ws.Range("F11").Value = 45.1384895657656
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 27.57448338579875
$ws.Range("J11").Value = 10.27649630325415
$ws.Range("N11").Value = 18.20245691061822

# Row 12
$ws.Range("B12").Value = 22.97227563811044
$ws.Range("C12").Value = 19.77578244515145
$ws.Range("D12").Value = 6.999919545511658
$ws.Range("E12").Value = 12.90482775805093
$ws.Range("F12").Value = 45.22824295537393
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 27.58543571784476
$ws.Range("J12").Value = 10.27756543703243
$ws.Range("N12").Value = 18.18193310800984

# Row 13
$ws.Range("B13").Value = 22.93932514299526
$ws.Range("C13").Value = 19.7437985761046
$ws.Range("D13").Value = 6.999196242228797
$ws.Range("E13").Value = 12.90182239162301
$ws.Range("F13").Value = 45.20883211338917
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 27.58302337128523
$ws.Range("J13").Value = 10.27732209588494
$ws.Range("N13").Value = 18.18634032150382

# Row 14
$ws.Range("B14").Value = 22.83167453372325
$ws.Range("C14").Value = 19.63926180056914
$ws.Range("D14").Value = 6.996853982151721
$ws.Range("E14").Value = 12.89208499392515
$ws.Range("F14").Value = 45.14583712083372
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 27.57536030696729
$ws.Range("J14").Value = 10.27657836557024
$ws.Range("N14").Value = 18.20076256490908

# Row 15
$ws.Range("B15").Value = 22.7655906645704
$ws.Range("C15").Value = 19.57505546142169
$ws.Range("D15").Value = 6.995432057065141
$ws.Range("E15").Value = 12.88616966825974
$ws.Range("F15").Value = 45.10748853077071
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 27.57082325818227
$ws.Range("J15").Value = 10.27616112121837
$ws.Range("N15").Value = 18.20963456760095

# Row 16
$ws.Range("B16").Value = 22.38478517165733
$ws.Range("C16").Value = 19.2045463652182
$ws.Range("D16").Value = 6.987480617656491
$ws.Range("E16").Value = 12.85302895870345
$ws.Range("F16").Value = 44.89141162515047
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 27.54722419169356
$ws.Range("J16").Value = 10.27435725342393
$ws.Range("N16").Value = 18.2610595816475

# Row 17
$ws.Range("B17").Value = 22.14942132504147
$ws.Range("C17").Value = 18.97507943351218
$ws.Range("D17").Value = 6.982779781204836
$ws.Range("E17").Value = 12.83337941316709
$ws.Range("F17").Value = 44.76218521014317
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 27.53488635312799
$ws.Range("J17").Value = 10.27377452511745
$ws.Range("N17").Value = 18.29312893217548

# Row 18
$ws.Range("B18").Value = 22.01345421336623
$ws.Range("C18").Value = 18.84234756145426
$ws.Range("D18").Value = 6.980141579060541
$ws.Range("E18").Value = 12.82232983657543
$ws.Range("F18").Value = 44.68909783332883
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 27.52857485462494
$ws.Range("J18").Value = 10.27363183290493
$ws.Range("N18").Value = 18.31176668390267

# Row 19
$ws.Range("B19").Value = 21.96732145644083
$ws.Range("C19").Value = 18.79728312698255
$ws.Range("D19").Value = 6.979259642412799
$ws.Range("E19").Value = 12.81863215818513
$ws.Range("F19").Value = 44.66456623745193
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 27.5265725560167
$ws.Range("J19").Value = 10.27361657158666
$ws.Range("N19").Value = 18.31811018895236

# Row 20
$ws.Range("B20").Value = 22.17453867652
$ws.Range("C20").Value = 18.99958519880242
$ws.Range("D20").Value = 6.983273415765604
$ws.Range("E20").Value = 12.83544507004975
$ws.Range("F20").Value = 44.77581358758945
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 27.53611847040658
$ws.Range("J20").Value = 10.27381663542519
$ws.Range("N20").Value = 18.2896952044285

# Row 21
$ws.Range("B21").Value = 22.86332456185387
$ws.Range("C21").Value = 19.67000341083666
$ws.Range("D21").Value = 6.997539300935979
$ws.Range("E21").Value = 12.89493488074392
$ws.Range("F21").Value = 45.16429084743028
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 27.57757845263361
$ws.Range("J21").Value = 10.27678883301286
$ws.Range("N21").Value = 18.19651849422862

# Row 22
$ws.Range("B22").Value = 23.30714356487744
$ws.Range("C22").Value = 20.10047833547329
$ws.Range("D22").Value = 7.007436115206691
$ws.Range("E22").Value = 12.93601847040728
$ws.Range("F22").Value = 45.42886583607406
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 27.61168959050491
$ws.Range("J22").Value = 10.28044598763717
$ws.Range("N22").Value = 18.13732273180366

# Row 23
$ws.Range("B23").Value = 23.07089865601851
$ws.Range("C23").Value = 19.87147501088029
$ws.Range("D23").Value = 7.002102085688088
$ws.Range("E23").Value = 12.91389195877141
$ws.Range("F23").Value = 45.2866980005666
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 27.59284102768794
$ws.Range("J23").Value = 10.27833720346477
$ws.Range("N23").Value = 18.16876159255137

# Row 24
$ws.Range("B24").Value = 22.16318514459021
$ws.Range("C24").Value = 18.98850864825006
$ws.Range("D24").Value = 6.983050042958884
$ws.Range("E24").Value = 12.83451041610467
$ws.Range("F24").Value = 44.76964843591502
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 27.53555899571116
$ws.Range("J24").Value = 10.27379699831368
$ws.Range("N24").Value = 18.2912469676159

# Row 25
$ws.Range("B25").Value = 21.16218809327524
$ws.Range("C25").Value = 18.00821276725975
$ws.Range("D25").Value = 6.965012919905337
$ws.Range("E25").Value = 12.75855131019764
$ws.Range("F25").Value = 44.25951293928659
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 27.50375114788037
$ws.Range("J25").Value = 10.27625639133188
$ws.Range("N25").Value = 18.43072715036169
